$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B text updates (wrapped as <token> strings)
$ws.Range("B2").Value  = "<there>"
$ws.Range("B4").Value  = "<it>"
$ws.Range("B5").Value  = "<has>"
$ws.Range("B6").Value  = "<coup>"
$ws.Range("B7").Value  = "<of>"
$ws.Range("B8").Value  = "<was>"
$ws.Range("B9").Value  = "<word>"
$ws.Range("B10").Value = "<paste>"
$ws.Range("B11").Value = "<which>"
$ws.Range("B12").Value = "<into>"
$ws.Range("B13").Value = "<ould>"
$ws.Range("B14").Value = "<many>"
$ws.Range("B15").Value = "<more>"

# Column C numeric updates
$ws.Range("C2").Value  = 51
$ws.Range("C3").Value  = 52
$ws.Range("C4").Value  = 54
$ws.Range("C5").Value  = 52
$ws.Range("C6").Value  = 50
$ws.Range("C7").Value  = 46
$ws.Range("C8").Value  = 52
$ws.Range("C9").Value  = 51
$ws.Range("C10").Value = 54
$ws.Range("C11").Value = 48
$ws.Range("C12").Value = 54
$ws.Range("C13").Value = 59
$ws.Range("C14").Value = 52
$ws.Range("C15").Value = 49
$ws.Range("C16").Value = 29
